# Auto-generated edit script: update Leve price/profit figures across 8 job sheets
# (data refresh from market board pricing; no formulas involved, values are static)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6: Days of Chunder / Antidote
$ws.Range("H6").Value = 2001.8334
$ws.Range("I6").Value = 2302.2
$ws.Range("J6").Value = 500
$ws.Range("K6").Value = 6906.599999999999
$ws.Range("L6").Value = 1500
$ws.Range("M6").Value = -6794.599999999999
$ws.Range("N6").Value = -1724

# Row 9: Distill, My Heart / Distilled Water
$ws.Range("H9").Value = 110.46667
$ws.Range("I9").Value = 77.90909000000001
$ws.Range("J9").Value = 200
$ws.Range("K9").Value = 77.90909000000001
$ws.Range("L9").Value = 200
$ws.Range("M9").Value = 91.09090999999999
$ws.Range("N9").Value = -538

# Row 12: Don't Be So Tallow / Beeswax
$ws.Range("H12").Value = 299
$ws.Range("I12").Value = 331.33334
$ws.Range("J12").Value = 202
$ws.Range("K12").Value = 331.33334
$ws.Range("L12").Value = 202
$ws.Range("M12").Value = -161.33334
$ws.Range("N12").Value = -542

# Row 16: Using Your Arcane Powers for Fun and Profit / Ash Picatrix
$ws.Range("H16").Value = 601.3333
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()

# Row 38: Just Give Him a Serum / Hi-Potion of Strength
$ws.Range("H38").Value = 2161.1667
$ws.Range("I38").Value = 1776.2858
$ws.Range("J38").Value = 2700
$ws.Range("K38").Value = 5328.857400000001
$ws.Range("L38").Value = 8100
$ws.Range("M38").Value = -4956.857400000001
$ws.Range("N38").Value = -8844

# Row 62: The Mustache Suits Him / Enchanted Mythrite Ink
$ws.Range("H62").Value = 4294.6313
$ws.Range("I62").Value = 2669.2307
$ws.Range("J62").Value = 7816.3335
$ws.Range("K62").Value = 2669.2307
$ws.Range("L62").Value = 7816.3335
$ws.Range("M62").Value = -2045.2307
$ws.Range("N62").Value = -9064.333500000001

# Row 65: Forgery of Convenience (L) / Enchanted Mythrite Ink
$ws.Range("H65").Value = 4294.6313
$ws.Range("I65").Value = 2669.2307
$ws.Range("J65").Value = 7816.3335
$ws.Range("K65").Value = 13346.1535
$ws.Range("L65").Value = 39081.6675
$ws.Range("M65").Value = -10226.1535
$ws.Range("N65").Value = -45321.6675

# Row 100: Asking for a Friend / Beetle Glue
$ws.Range("H100").Value = 47621650
$ws.Range("I100").Value = 33334974
$ws.Range("J100").Value = 83338340
$ws.Range("K100").Value = 33334974
$ws.Range("L100").Value = 83338340
$ws.Range("M100").Value = -33334433
$ws.Range("N100").Value = -83339422

# Row 111: An Eye for Healing / Grade 1 Dexterity Alkahest
$ws.Range("H111").Value = 3000
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 3000
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 9000
$ws.Range("M111").ClearContents()
$ws.Range("N111").Value = -15134

# Row 113: Amaro Kart / Starch Glue
$ws.Range("H113").Value = 252126.25
$ws.Range("I113").Value = 501252.5
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 501252.5
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = -497998.5
$ws.Range("N113").Value = -9508

# Row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 267082.5
$ws.Range("I132").Value = 338224.22
$ws.Range("J132").Value = 34255
$ws.Range("K132").Value = 1014672.66
$ws.Range("L132").Value = 102765
$ws.Range("M132").Value = -1012142.66
$ws.Range("N132").Value = -107825

# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 1547.6
$ws.Range("I137").Value = 912.5
$ws.Range("J137").Value = 1846.4706
$ws.Range("K137").Value = 2737.5
$ws.Range("L137").Value = 5539.4118
$ws.Range("M137").Value = -187.5
$ws.Range("N137").Value = -10639.4118

# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 6561091
$ws.Range("I138").Value = 1009983.8
$ws.Range("J138").Value = 10002778
$ws.Range("K138").Value = 3029951.4
$ws.Range("L138").Value = 30008334
$ws.Range("M138").Value = -3024811.4
$ws.Range("N138").Value = -30018614

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 3853.8096
$ws.Range("I32").Value = 3309.5186
$ws.Range("J32").Value = 4833.533
$ws.Range("K32").Value = 3309.5186
$ws.Range("L32").Value = 4833.533
$ws.Range("M32").Value = -3022.5186
$ws.Range("N32").Value = -5407.533

# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 5063.0625
$ws.Range("I74").Value = 1427.85
$ws.Range("J74").Value = 11121.75
$ws.Range("K74").Value = 1427.85
$ws.Range("L74").Value = 11121.75
$ws.Range("M74").Value = -553.8499999999999
$ws.Range("N74").Value = -12869.75

# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 5063.0625
$ws.Range("I77").Value = 1427.85
$ws.Range("J77").Value = 11121.75
$ws.Range("K77").Value = 7139.25
$ws.Range("L77").Value = 55608.75
$ws.Range("M77").Value = -2771.25
$ws.Range("N77").Value = -64344.75

# Row 86: Sir, Dost Thou Even Heft / Adamantite Chain Hose of Fending
$ws.Range("H86").Value = 50000
$ws.Range("J86").Value = 50000
$ws.Range("L86").Value = 50000
$ws.Range("N86").Value = -52372

# Row 89: Men in Adamantite (L) / Adamantite Chain Hose of Fending
$ws.Range("H89").Value = 50000
$ws.Range("J89").Value = 50000
$ws.Range("L89").Value = 150000
$ws.Range("N89").Value = -161856

$ws = $wb.Worksheets.Item("BSM")
# Row 105: Ingot to Wing It / Molybdenum Ingot
$ws.Range("H105").Value = 17546124
$ws.Range("I105").Value = 19610138
$ws.Range("J105").Value = 2005.5
$ws.Range("K105").Value = 19610138
$ws.Range("L105").Value = 2005.5
$ws.Range("M105").Value = -19608391
$ws.Range("N105").Value = -5499.5

# Row 139: Maul Me / Titanium Gold Maul
$ws.Range("H139").Value = 95000
$ws.Range("J139").Value = 95000
$ws.Range("L139").Value = 95000
$ws.Range("N139").Value = -105280

$ws = $wb.Worksheets.Item("CRP")
# Row 20: Re-crating the Scene / Iron Spear
$ws.Range("H20").Value = 49832.668
$ws.Range("J20").Value = 49832.668
$ws.Range("L20").Value = 49832.668
$ws.Range("N20").Value = -50304.668

# Row 30: Polearms Aplenty / Iron Spear
$ws.Range("H30").Value = 49832.668
$ws.Range("J30").Value = 49832.668
$ws.Range("L30").Value = 49832.668
$ws.Range("N30").Value = -50014.668

# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 2206.7827
$ws.Range("I31").Value = 1465
$ws.Range("J31").Value = 3902.2856
$ws.Range("K31").Value = 1465
$ws.Range("L31").Value = 3902.2856
$ws.Range("M31").Value = -1170
$ws.Range("N31").Value = -4492.2856

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 2206.7827
$ws.Range("I34").Value = 1465
$ws.Range("J34").Value = 3902.2856
$ws.Range("K34").Value = 1465
$ws.Range("L34").Value = 3902.2856
$ws.Range("M34").Value = -1263
$ws.Range("N34").Value = -4306.2856

# Row 128: An A-prop-riate Request / Ironwood Spear
$ws.Range("H128").Value = 49832.668
$ws.Range("J128").Value = 49832.668
$ws.Range("L128").Value = 49832.668
$ws.Range("N128").Value = -59792.668

# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 3045.3809
$ws.Range("I132").Value = 1929.5
$ws.Range("J132").Value = 5277.143
$ws.Range("K132").Value = 5788.5
$ws.Range("L132").Value = 15831.429
$ws.Range("M132").Value = -3258.5
$ws.Range("N132").Value = -20891.429

$ws = $wb.Worksheets.Item("CUL")
# Row 76: Old Victories, New Tastes / Dhalmel Fricassee
$ws.Range("H76").Value = 5000
$ws.Range("J76").Value = 5000
$ws.Range("L76").Value = 15000
$ws.Range("N76").Value = -15766

# Row 79: The Eats of Authenticity (L) / Dhalmel Fricassee
$ws.Range("H79").Value = 5000
$ws.Range("J79").Value = 5000
$ws.Range("L79").Value = 15000
$ws.Range("N79").Value = -17652

# Row 107: Slippery Service / Frantoio Oil
$ws.Range("H107").Value = 1286.7258
$ws.Range("I107").Value = 378.86667
$ws.Range("J107").Value = 1576.4681
$ws.Range("K107").Value = 1136.60001
$ws.Range("L107").Value = 4729.4043
$ws.Range("M107").Value = 783.3999899999999
$ws.Range("N107").Value = -8569.4043

# Row 113: Can't Eat Just One / Night Vinegar
$ws.Range("H113").Value = 22727710
$ws.Range("I113").Value = 810
$ws.Range("J113").Value = 23809944
$ws.Range("K113").Value = 2430
$ws.Range("L113").Value = 71429832
$ws.Range("M113").Value = -260
$ws.Range("N113").Value = -71434172

# Row 127: A Stickler for Carrots / Carrot Nibbles
$ws.Range("H127").Value = 3187.6875
$ws.Range("J127").Value = 3187.6875
$ws.Range("L127").Value = 9563.0625
$ws.Range("N127").Value = -19483.0625

# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 3726.4565
$ws.Range("I131").Value = 765
$ws.Range("J131").Value = 4008.5
$ws.Range("K131").Value = 2295
$ws.Range("L131").Value = 12025.5
$ws.Range("M131").Value = 2745
$ws.Range("N131").Value = -22105.5

$ws = $wb.Worksheets.Item("GSM")
# Row 122: Awarding Academic Excellence / Ametrine
$ws.Range("H122").Value = 2780404.5
$ws.Range("I122").Value = 2780404.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8341213.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -8338763.5
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban / Leather
$ws.Range("H7").Value = 2596.6538
$ws.Range("I7").Value = 1336
$ws.Range("J7").Value = 2761.087
$ws.Range("K7").Value = 1336
$ws.Range("L7").Value = 2761.087
$ws.Range("M7").Value = -1224
$ws.Range("N7").Value = -2985.087

# Row 40: Best Served Toad / Toad Leather
$ws.Range("H40").Value = 4361.1113
$ws.Range("I40").Value = 2500
$ws.Range("K40").Value = 2500
$ws.Range("M40").Value = -2364

# Row 61: Spelling Me Softly / Raptor Leather
$ws.Range("H61").Value = 3573
$ws.Range("I61").Value = 3992.2942
$ws.Range("K61").Value = 3992.2942
$ws.Range("M61").Value = -3790.2942

# Row 113: Peace in Rest / Atrociraptor Leather
$ws.Range("H113").Value = 3573
$ws.Range("I113").Value = 3992.2942
$ws.Range("K113").Value = 3992.2942
$ws.Range("M113").Value = -1822.2942

# Row 122: Hell on Leather / Gaja Leather
$ws.Range("H122").Value = 3364.9
$ws.Range("I122").Value = 2562.375
$ws.Range("J122").Value = 3899.9167
$ws.Range("K122").Value = 7687.125
$ws.Range("L122").Value = 11699.7501
$ws.Range("M122").Value = -5237.125
$ws.Range("N122").Value = -16599.7501

# Row 126: Battered Books / Saiga Leather
$ws.Range("H126").Value = 2596.6538
$ws.Range("I126").Value = 1336
$ws.Range("J126").Value = 2761.087
$ws.Range("K126").Value = 4008
$ws.Range("L126").Value = 8283.261
$ws.Range("M126").Value = -1538
$ws.Range("N126").Value = -13223.261

# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 3474.8518
$ws.Range("I132").Value = 1999.7142
$ws.Range("J132").Value = 5063.4614
$ws.Range("K132").Value = 5999.142599999999
$ws.Range("L132").Value = 15190.3842
$ws.Range("M132").Value = -3469.142599999999
$ws.Range("N132").Value = -20250.3842

$ws = $wb.Worksheets.Item("WVR")
# Row 64: Ribbon of Remembrance / Rainbow Ribbon of Healing
$ws.Range("H64").Value = 30000
$ws.Range("J64").Value = 30000
$ws.Range("L64").Value = 30000
$ws.Range("N64").Value = -30496

# Row 67: The Road Was a Ribbon of Moonlight (L) / Rainbow Ribbon of Healing
$ws.Range("H67").Value = 30000
$ws.Range("J67").Value = 30000
$ws.Range("L67").Value = 30000
$ws.Range("N67").Value = -31716

# Row 113: A Tender Table / Pixie Floss
$ws.Range("H113").Value = 839.3889
$ws.Range("I113").Value = 1019.9286
$ws.Range("J113").Value = 207.5
$ws.Range("K113").Value = 3059.7858
$ws.Range("L113").Value = 622.5
$ws.Range("M113").Value = -889.7857999999997
$ws.Range("N113").Value = -4962.5

# Row 123: Helping Handwear / Fingerless Darkhempen Gloves of Healing
$ws.Range("H123").Value = 33280.668
$ws.Range("J123").Value = 33280.668
$ws.Range("L123").Value = 33280.668
$ws.Range("N123").Value = -43080.668

# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 15627135
$ws.Range("I132").Value = 25001772
$ws.Range("J132").Value = 2741.6667
$ws.Range("K132").Value = 75005316
$ws.Range("L132").Value = 8225.000100000001
$ws.Range("M132").Value = -75002786
$ws.Range("N132").Value = -13285.0001

# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 11495550
$ws.Range("I136").Value = 17544508
$ws.Range("J136").Value = 2530.3
$ws.Range("K136").Value = 52633524
$ws.Range("L136").Value = 7590.900000000001
$ws.Range("M136").Value = -52630974
$ws.Range("N136").Value = -12690.9
